$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C47").Value2 = "Test Paper Name"
$ws.Rows.Item(47).RowHeight = 210
Write-Host "RowHeight47: $($ws.Rows.Item(47).RowHeight)"
Write-Host "C47: $($ws.Range('C47').Value2)"

$ws.Range("K46").Font.Bold = $true
Write-Host "K46 bold: $($ws.Range('K46').Font.Bold)"

$r = $ws.Range("A47")
Write-Host "A47 value: $($r.Value2)"
